$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new rows at 761 (pushes old rows 761-818 down to 764-821,
# matching the dimension growing from A1:T818 to A1:T821).
$ws.Rows("761:763").Insert()

# Fill the 3 new rows with the new weekly record (date 44714 = 2022-06-02,
# Hass avocado, Especial/Primera/Segunda quality, Provincia de Limari origin).
$newRows = @(
    @{ Row = 761; L = "Especial"; M = 300; N = 3300; O = 3400; P = 3350; S = 3350 },
    @{ Row = 762; L = "Primera";  M = 300; N = 3100; O = 3200; P = 3150; S = 3150 },
    @{ Row = 763; L = "Segunda";  M = 240; N = 2900; O = 3000; P = 2950; S = 2950 }
)

foreach ($r in $newRows) {
    $row = $r.Row
    $ws.Cells.Item($row, 1).Value = 8
    $ws.Cells.Item($row, 2).Value = "Terminal La Palmera de La Serena"
    $ws.Cells.Item($row, 3).Value = "Coquimbo"
    $ws.Cells.Item($row, 4).Value = 44714
    $ws.Cells.Item($row, 5).Value = 4
    $ws.Cells.Item($row, 6).Value = "Fruta"
    $ws.Cells.Item($row, 7).Value = 100106
    $ws.Cells.Item($row, 8).Value = "Oleaginosos"
    $ws.Cells.Item($row, 9).Value = 100106002
    $ws.Cells.Item($row, 10).Value = "Palta"
    $ws.Cells.Item($row, 11).Value = "Hass"
    $ws.Cells.Item($row, 12).Value = $r.L
    $ws.Cells.Item($row, 13).Value = $r.M
    $ws.Cells.Item($row, 14).Value = $r.N
    $ws.Cells.Item($row, 15).Value = $r.O
    $ws.Cells.Item($row, 16).Value = $r.P
    $ws.Cells.Item($row, 17).Value = "`$/kilo (en caja de 17 kilos)"
    $ws.Cells.Item($row, 18).Value = "Provincia de Limar$([char]0x00ED)"
    $ws.Cells.Item($row, 19).Value = $r.S
    $ws.Cells.Item($row, 20).Value = 1
}
